# Updated symbol list on Mon Dec 26 16:11:12 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Hora" (column G) values in the
# crypto symbol table on Sheet1 (data rows 2-51). Both columns store their
# values as plain text, so every new value below is written with a leading
# apostrophe - this keeps Excel from auto-converting the numeric-looking
# text into a real number, matching how the cells were already stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates -------------------------------------------
$ws.Range("D2").Value = "'242.47"
$ws.Range("D3").Value = "'23.02"
$ws.Range("D4").Value = "'5.430"
$ws.Range("D5").Value = "'0.05891"
$ws.Range("D6").Value = "'3.435"
$ws.Range("D7").Value = "'6.547"
$ws.Range("D8").Value = "'0.8108"
$ws.Range("D9").Value = "'0.9363"
$ws.Range("D10").Value = "'0.1414"
$ws.Range("D11").Value = "'0.07439"
$ws.Range("D12").Value = "'0.03268"
$ws.Range("D13").Value = "'0.03060"
$ws.Range("D14").Value = "'0.09334"
$ws.Range("D15").Value = "'3.850"
$ws.Range("D16").Value = "'0.001573"
$ws.Range("D18").Value = "'0.0005954"
$ws.Range("D19").Value = "'0.005854"
$ws.Range("D20").Value = "'0.001251"
$ws.Range("D21").Value = "'0.004894"
$ws.Range("D22").Value = "'0.00006805"
$ws.Range("D23").Value = "'3.589"
$ws.Range("D24").Value = "'2.126"
$ws.Range("D25").Value = "'0.3230"
$ws.Range("D26").Value = "'0.1308"
$ws.Range("D27").Value = "'0.0002286"
$ws.Range("D40").Value = "'0.03926"
$ws.Range("D42").Value = "'0.1070"
$ws.Range("D43").Value = "'0.002562"
$ws.Range("D44").Value = "'0.009255"
$ws.Range("D45").Value = "'0.00005202"
$ws.Range("D47").Value = "'0.7305"
$ws.Range("D48").Value = "'0.002369"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.0002001"

# --- Hora (column G) updates: every data row moves from 15 to 16 --------
$ws.Range("G2").Value = "'16"
$ws.Range("G3").Value = "'16"
$ws.Range("G4").Value = "'16"
$ws.Range("G5").Value = "'16"
$ws.Range("G6").Value = "'16"
$ws.Range("G7").Value = "'16"
$ws.Range("G8").Value = "'16"
$ws.Range("G9").Value = "'16"
$ws.Range("G10").Value = "'16"
$ws.Range("G11").Value = "'16"
$ws.Range("G12").Value = "'16"
$ws.Range("G13").Value = "'16"
$ws.Range("G14").Value = "'16"
$ws.Range("G15").Value = "'16"
$ws.Range("G16").Value = "'16"
$ws.Range("G17").Value = "'16"
$ws.Range("G18").Value = "'16"
$ws.Range("G19").Value = "'16"
$ws.Range("G20").Value = "'16"
$ws.Range("G21").Value = "'16"
$ws.Range("G22").Value = "'16"
$ws.Range("G23").Value = "'16"
$ws.Range("G24").Value = "'16"
$ws.Range("G25").Value = "'16"
$ws.Range("G26").Value = "'16"
$ws.Range("G27").Value = "'16"
$ws.Range("G28").Value = "'16"
$ws.Range("G29").Value = "'16"
$ws.Range("G30").Value = "'16"
$ws.Range("G31").Value = "'16"
$ws.Range("G32").Value = "'16"
$ws.Range("G33").Value = "'16"
$ws.Range("G34").Value = "'16"
$ws.Range("G35").Value = "'16"
$ws.Range("G36").Value = "'16"
$ws.Range("G37").Value = "'16"
$ws.Range("G38").Value = "'16"
$ws.Range("G39").Value = "'16"
$ws.Range("G40").Value = "'16"
$ws.Range("G41").Value = "'16"
$ws.Range("G42").Value = "'16"
$ws.Range("G43").Value = "'16"
$ws.Range("G44").Value = "'16"
$ws.Range("G45").Value = "'16"
$ws.Range("G46").Value = "'16"
$ws.Range("G47").Value = "'16"
$ws.Range("G48").Value = "'16"
$ws.Range("G49").Value = "'16"
$ws.Range("G50").Value = "'16"
$ws.Range("G51").Value = "'16"

